# Update all Logic Component Related Diagrams and ppt files
#
# 1) Refresh the "datetimeFigureOut" date placeholders (slide master,
#    every slide layout, and the notes master) from 7/11/2017 -> 4/14/2018.
# 2) Rename the ":AddressBookParser" lifeline box to ":CatalogueParser"
#    (now written as a single centred paragraph) and shrink/shift it.
# 3) Rename the deletePerson(p) call label to deleteBook(p).

$p = $ppt.ActivePresentation
$newDate = "4/14/2018"

function Set-DatePlaceholderText {
    param($shapes, [string]$text)

    foreach ($sh in $shapes) {
        $isDatePlaceholder = $false
        try {
            if (($sh.Type -eq 14) -and ($sh.PlaceholderFormat.Type -eq 16)) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

# --- Slide master date placeholder ---
$slideMaster = $p.SlideMaster
Set-DatePlaceholderText $slideMaster.Shapes $newDate

# --- Every slide layout's date placeholder ---
$layouts = $slideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# --- Notes master date placeholder ---
$notesMaster = $p.NotesMaster
Set-DatePlaceholderText $notesMaster.Shapes $newDate

# --- Slide 1 shape edits ---
$slide = $p.Slides.Item(1)

foreach ($shape in $slide.Shapes) {

    if ($shape.Id -eq 16) {
        # ":AddressBookParser" (two paragraphs) -> ":CatalogueParser" (one paragraph)
        $shape.Left = (2784016 + 0.5) / 12700.0
        $shape.Width = (1104348 + 0.5) / 12700.0

        $tr = $shape.TextFrame.TextRange

        # Drop "Address", leaving ":" + CRLF + "BookParser"
        $tr.Characters(2, 7).Text = ""

        # "Book" -> "Catalogue" inside the second run, preserving its own formatting
        $tr2 = $shape.TextFrame.TextRange
        $tr2.Characters(3, 4).Text = "Catalogue"

        # Merge the two paragraphs into a single centred paragraph made of two runs
        $finalText = $shape.TextFrame.TextRange.Text -replace "[\r\n\v]", ""
        $shape.TextFrame.TextRange.Text = $finalText
        $merged = $shape.TextFrame.TextRange
        $merged.Characters(2, $finalText.Length - 1).Font.Size = 16
    }

    if ($shape.Id -eq 78) {
        # deletePerson(p) -> deleteBook(p), keep the "(p)" run untouched
        $shape.TextFrame.TextRange.Characters(1, 12).Text = "deleteBook"
    }
}
